$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C contains a "last changed" date (serial 45189 = 2023-09-20).
# Update every populated row (2 through 260) to the new date serial 45190
# (2023-09-21), matching the diff which bumps every C-column cell by +1.
$lastRow = 260
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45190
}
